$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 78, shifting existing rows 78:171 down to 79:172
$ws.Rows("78:78").Insert()

# Populate the newly inserted row 78 with fresh data
$ws.Range("A78").Value = 10
$ws.Range("B78").Value = "Vega Modelo de Temuco"
$ws.Range("C78").Value = "La Araucanía"
$ws.Range("D78").Value = 44413
$ws.Range("E78").Value = 9
$ws.Range("F78").Value = "Fruta"
$ws.Range("G78").Value = 100108
$ws.Range("H78").Value = "Tropicales y subtropicales"
$ws.Range("I78").Value = 100108002
$ws.Range("J78").Value = "Mango"
$ws.Range("K78").Value = "Sin especificar"
$ws.Range("L78").Value = "Primera"
$ws.Range("M78").Value = 650
$ws.Range("N78").Value = 9000
$ws.Range("O78").Value = 9000
$ws.Range("P78").Value = 9000
$ws.Range("Q78").Value = '$/bandeja 4 kilos'
$ws.Range("R78").Value = "Brasil"
$ws.Range("S78").Value = 2250
$ws.Range("T78").Value = 4
